$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new task row (row 12) under the "Investigacion" section for the
#     UPnP research/attempt task, shifting every following row down by one. ---
$ws.Rows("12:12").Insert()

# The freshly inserted row inherits formatting (and an empty cell) from the row
# above; strip that stray cell before filling in the real content.
$ws.Range("AH12").Style = "Normal"
$ws.Range("AH12").ClearContents()

# Copy the task-name formatting (style used by other task rows) onto A12, and
# the hours-entry formatting (yellow fill, right aligned) onto BB12.
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("AH11").Copy()
$ws.Range("BB12").PasteSpecial(-4122)

# --- New hour entries added to existing rows (now shifted down by the insert
#     above): "Implementación notificaciones" row gained a BA cell, and the
#     renamed "selección personaje y mapa" task plus the row below it gained
#     new cells too. ---
$ws.Range("AU31").Copy()
$ws.Range("BA31").PasteSpecial(-4122)
$ws.Range("BA31").Value = "4 h."

$ws.Range("AZ34").Copy()
$ws.Range("BB34").PasteSpecial(-4122)
$ws.Range("BB34").Value = "2.5 h."

$ws.Range("AZ35").Copy()
$ws.Range("BA35").PasteSpecial(-4122)
$ws.Range("BA35").Value = "2.5 h."

# --- Rename the task to mention the new "final partida" (match end) screen. ---
$ws.Range("A34").Value = "Implementación selección personaje y mapa, y final partida"

# --- Fill in the new UPnP research task (after the rename above, so the
#     shared-string table ends up in the same order as the authored edit). ---
$ws.Range("A12").Value = "UPnP (e intento de implementación)"
$ws.Range("BB12").Value = "3.5 h."

# --- Cosmetic sheet-level tweaks from the diff. ---
$ws.Columns("A:A").ColumnWidth = 59.92

[void]$ws.Range("BB12").Select()
